# Metodos para usar tecla enter en el pryecto
# Fill in the "Progresos" (D column) values that were pending (0) and
# bump the C3 value, simulating typing values and pressing Enter which
# moves the active cell down one row at a time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 was left at 0 - fill it in along the way
$ws.Range("C3").Value = 70

# Column D ("progreso" column) values entered one by one (Enter moves
# the selection down after each entry)
$ws.Range("D2").Value = 100
$ws.Range("D3").Value = 80
$ws.Range("D4").Value = 88
$ws.Range("D5").Value = 60
$ws.Range("D6").Value = 60
$ws.Range("D7").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 85
$ws.Range("D10").Value = 85
$ws.Range("D11").Value = 96
$ws.Range("D12").Value = 80
$ws.Range("D13").Value = 90
$ws.Range("D14").Value = 100

# After typing the last value and pressing Enter, the selection moves
# down to D15.
$ws.Range("D15").Select()
